$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 206, shifting rows 206:272 down to 207:273.
$ws.Rows(206).Insert()

# The values that used to be in row 206 are now in row 207; copy them into
# the new row 206, then overwrite the Fecha (D) and Volumen (J) values.
for ($c = 1; $c -le 18; $c++) {
    $src = $ws.Cells.Item(207, $c)
    $dst = $ws.Cells.Item(206, $c)
    $dst.Value2 = $src.Value2
}

$ws.Cells.Item(206, 4).Value2 = 44524
$ws.Cells.Item(206, 10).Value2 = 500
